$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 6;   D = "[1, 0, 0, 0, 0, 0, 0]"; E = "['Normal']" },
    @{ Row = 8;   D = "[1, 0, 1, 0, 0, 0, 0]"; E = "['Normal', 'HardwareFault']" },
    @{ Row = 11;  D = "[1, 0, 1, 0, 1, 0, 0]"; E = "['Normal', 'HardwareFault', 'RegulationViolation']" },
    @{ Row = 16;  D = "[1, 0, 0, 0, 1, 0, 0]"; E = "['Normal', 'RegulationViolation']" },
    @{ Row = 24;  D = "[1, 0, 0, 0, 0, 0, 0]"; E = "['Normal']" },
    @{ Row = 27;  D = "[1, 0, 0, 0, 0, 0, 1]"; E = "['Normal', 'SoftwareFault']" },
    @{ Row = 28;  D = "[1, 0, 0, 0, 0, 0, 1]"; E = "['Normal', 'SoftwareFault']" },
    @{ Row = 29;  D = "[1, 0, 0, 0, 0, 0, 1]"; E = "['Normal', 'SoftwareFault']" },
    @{ Row = 36;  D = "[1, 0, 1, 0, 0, 0, 0]"; E = "['Normal', 'HardwareFault']" },
    @{ Row = 38;  D = "[1, 0, 0, 0, 0, 0, 0]"; E = "['Normal']" },
    @{ Row = 56;  D = "[1, 0, 0, 0, 0, 0, 0]"; E = "['Normal']" },
    @{ Row = 58;  D = "[1, 0, 0, 0, 0, 0, 0]"; E = "['Normal']" },
    @{ Row = 61;  D = "[1, 0, 0, 0, 0, 0, 0]"; E = "['Normal']" },
    @{ Row = 80;  D = "[1, 0, 0, 0, 0, 0, 0]"; E = "['Normal']" },
    @{ Row = 83;  D = "[1, 0, 0, 0, 0, 0, 0]"; E = "['Normal']" },
    @{ Row = 107; D = "[1, 0, 0, 0, 0, 1, 0]"; E = "['Normal', 'CommunicationIssue']" },
    @{ Row = 109; D = "[1, 0, 0, 0, 0, 0, 1]"; E = "['Normal', 'SoftwareFault']" },
    @{ Row = 113; D = "[1, 0, 1, 0, 0, 0, 0]"; E = "['Normal', 'HardwareFault']" }
)

foreach ($u in $updates) {
    $ws.Range("D" + $u.Row).Value = $u.D
    $ws.Range("E" + $u.Row).Value = $u.E
}
